$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D2 gains a real value: unit "fL" (平均红细胞体积 row now has a unit too)
$ws.Range("D2").Value = "fL"

# D4, D5, D6, D10 gain an (empty) cell in the "单位" column, mirroring the
# other rows of the table that already carry a D-column cell. Re-apply the
# default "Normal" style to force Excel to materialize the blank cell
# without actually changing its value/format.
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D10").Style = "Normal"
